# AFDP-9108 Fix Consultation Module Admin Issues
# - Insert default type on Consultation save: add a "Set Consultation Type"
#   rule row into the "Save Consultation Rules" rule table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The rule table body (rows 22-33) uses style 17 for columns B/C/D and
# style 1 for column A. Insert a brand-new row above the existing
# "Set Consultation Priority" row (row 25) so every row below shifts down
# by one, then copy the formatting of the row that is now directly below
# the insertion point (the old "Set Consultation Priority" row, now at 26)
# into the new blank row so it keeps the correct cell styles/borders.
$ws.Rows.Item(25).Insert()
$ws.Range("B26:D26").Copy($ws.Range("B25:D25"))

# Populate the new rule row with the "Set Consultation Type" rule.
# Values are written in the same order the original workbook introduced
# the corresponding shared strings (Rule Name, Action, then Condition).
$ws.Range("B25").Value = "Set Consultation Type"
$ws.Range("D25").Value = "setConsultationType, 'Consultation'"
$ws.Range("C25").Value = 'consultationType == null || consultationType.equals("")'

# Restore/update the active selection like the saved workbook shows.
$ws.Range("C28").Select()
